# TMTC0032668_VerifyActivityIsLinkedToTheRelatedCampaign.xlsx
# Test data update (8th July 2025): the "Campaign" sheet's RecordType value
# (cell E2) changes from "BAS" to "None".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Campaign")
$ws.Range("E2").Value = "None"
